# "Generate Report for Archive" — refresh localization status cells that were
# still showing the old "Ready for handoff" status to "In Translation", then
# re-size the affected Status columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $touchedCols = New-Object System.Collections.Generic.HashSet[int]

    foreach ($cell in $used.Cells) {
        $val = [string]$cell.Value2
        if ($val -eq "Ready for handoff") {
            $cell.Value = "In Translation"
            [void]$touchedCols.Add($cell.Column)
        }
    }

    foreach ($colIndex in $touchedCols) {
        $col = $ws.Cells.Item(1, $colIndex).EntireColumn
        $col.AutoFit() | Out-Null
        $col.ColumnWidth = 12.5
    }
}
